$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.309.27'
$ws.Range('E2').Value = '  -1.28%  '
$ws.Range('D3').Value = '2.045.97'
$ws.Range('E3').Value = '  -1.48%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '228.75'
$ws.Range('E5').Value = '  -1.82%  '
$ws.Range('D6').Value = '0.613'
$ws.Range('E6').Value = '  -1.88%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '56.56'
$ws.Range('E8').Value = '  -3.35%  '
$ws.Range('D9').Value = '0.384'
$ws.Range('E9').Value = '  -2.28%  '
$ws.Range('D10').Value = '0.0785'
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('D12').Value = '14.75'
$ws.Range('E12').Value = '  -0.50%  '
$ws.Range('D13').Value = '2.333.30'
$ws.Range('E13').Value = '  -2.12%  '
$ws.Range('D14').Value = '20.67'
$ws.Range('E14').Value = '  -2.02%  '
$ws.Range('E15').Value = '  -3.50%  '
$ws.Range('D16').Value = '5.28'
$ws.Range('E16').Value = '  -1.46%  '
$ws.Range('D17').Value = '2.060.39'
$ws.Range('E17').Value = '  -1.03%  '
$ws.Range('D18').Value = '37.206.26'
$ws.Range('E18').Value = '  -1.30%  '
$ws.Range('D19').Value = '6.04'
$ws.Range('E19').Value = '  -1.36%  '
$ws.Range('D20').Value = "'69.30"
$ws.Range('E20').Value = '  -3.27%  '
$ws.Range('D21').Value = '0.0₃0825'
$ws.Range('E21').Value = '  -2.07%  '
$ws.Range('D22').Value = '225.68'
$ws.Range('E22').Value = '  -1.56%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').Value = '2.28'
$ws.Range('E25').Value = '  -5.10%  '
$ws.Range('D26').Value = '9.67'
$ws.Range('E26').Value = '  -0.57%  '
$ws.Range('D27').Value = '166.18'
$ws.Range('E27').Value = '  -3.29%  '
$ws.Range('E28').Value = '  -7.15%  '
$ws.Range('D29').Value = '18.98'
$ws.Range('E29').Value = '  -2.30%  '
$ws.Range('E30').Value = '  -3.83%  '
$ws.Range('E31').Value = '  -2.00%  '
$ws.Range('D32').Value = '4.51'
$ws.Range('E32').Value = '  -4.78%  '
$ws.Range('D33').Value = '0.0615'
$ws.Range('E33').Value = '  -2.86%  '
$ws.Range('D34').Value = '4.57'
$ws.Range('E34').Value = '  -1.92%  '
$ws.Range('D35').Value = '2.44'
$ws.Range('E35').Value = '  -0.36%  '
$ws.Range('E36').Value = '  +1.16%  '
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('D38').Value = '3.24'
$ws.Range('E38').Value = '  -4.72%  '
$ws.Range('E39').Value = '  -3.27%  '
$ws.Range('E40').Value = '  -4.91%  '
$ws.Range('D41').Value = '1.482.17'
$ws.Range('E41').Value = '  +2.13%  '
$ws.Range('E42').Value = '  -0.85%  '
$ws.Range('D43').Value = '16.89'
$ws.Range('E43').Value = '  -0.72%  '
$ws.Range('D44').Value = '0.0941'
$ws.Range('E44').Value = '  -3.60%  '
$ws.Range('D45').Value = '96.35'
$ws.Range('E45').Value = '  -5.11%  '
$ws.Range('E46').Value = '  +0.59%  '
$ws.Range('D47').Value = '3.96'
$ws.Range('E47').Value = '  -3.29%  '
$ws.Range('E48').Value = '  -4.35%  '
$ws.Range('D49').Value = '7.12'
$ws.Range('E49').Value = '  -3.76%  '
$ws.Range('E50').Value = '  -2.37%  '
$ws.Range('D51').Value = '2.236.63'
$ws.Range('E51').Value = '  -1.47%  '
